$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.630.35"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "'1.848.01"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'312.51"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4283"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("D10").Value = "'0.07317"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").Value = "'0.8756"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").Value = "'20.69"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'1.847.86"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "'5.328"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "'6.520"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "'0.06911"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'79.99"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("D19").Value = "'0.000009010"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'15.33"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").Value = "'27.655.69"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("D25").Value = "'2.076.23"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "'1.988"
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("D27").Value = "'155.12"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").Value = "'18.78"
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("D29").Value = "'121.49"
$ws.Range("E29").Value = "  +8.76%  "
$ws.Range("D30").Value = "'5.290"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").Value = "'1.848"
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").Value = "'0.08903"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("D34").Value = "'4.553"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'2.967"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("D36").Value = "'1.103"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "'0.05411"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "'1.088"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").Value = "'0.01934"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "'2.816"
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("D41").Value = "'0.5078"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").Value = "'6.761"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "'8.372"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "'0.06550"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("D46").Value = "'10.37"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'105.19"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("D48").Value = "'0.4679"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'1.621"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "'64.36"
$ws.Range("E51").Value = "  -0.75%  "
